# Weekly data refresh: a new price-report row for "Orégano" at the
# Vega Central Mapocho de Santiago market is inserted at row 52 (the
# sheet is sorted by date, so the existing rows 52-100 shift down to
# 53-101). The new row carries the latest week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing row 52 (and everything below it) down by one row.
$ws.Rows.Item(52).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(52, 1).Value  = 9
$ws.Cells.Item(52, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(52, 3).Value  = "Metropolitana"
$ws.Cells.Item(52, 4).Value  = 45049
$ws.Cells.Item(52, 5).Value  = 13
$ws.Cells.Item(52, 6).Value  = 100112029
$ws.Cells.Item(52, 7).Value  = "Orégano"
$ws.Cells.Item(52, 8).Value  = "Sin especificar"
$ws.Cells.Item(52, 9).Value  = "Primera"
$ws.Cells.Item(52, 10).Value = 16
$ws.Cells.Item(52, 11).Value = 17000
$ws.Cells.Item(52, 12).Value = 17000
$ws.Cells.Item(52, 13).Value = 17000
$ws.Cells.Item(52, 14).Value = "$/docena de atados"
$ws.Cells.Item(52, 15).Value = "Región Metropolitana"
$ws.Cells.Item(52, 16).Value = 5667
$ws.Cells.Item(52, 17).Value = 3
$ws.Cells.Item(52, 18).Value = "Hortaliza"
